$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Maçlar")

# Fix typo in shared string: "Fortuna Unıted" -> "Fortuna United"
# This string appears in cell D3 / D8 (Takim1 column reference with shared string index 21)
foreach ($cell in $ws.UsedRange.Cells) {
    if ($cell.Value2 -eq "Fortuna Unıted") {
        $cell.Value2 = "Fortuna United"
    }
}

# Add the missing match result for row 9 (F9 = 1, G9 = 1)
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1

# Update the active selection to F10
$ws.Range("F10").Select()
